# Populate email-parsing table (rows 2-8) added on 06022020.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- raw source data (Email in column D, Password in column J) ---------
$ws.Range("D2").Value = "sourabh.awasthi@capgemini.com"
$ws.Range("J2").Value = "o8{5D}@3T7Dd"

$ws.Range("D3").Value = "sandipan.deb@capgemini.com"
$ws.Range("J3").Value = "p@#35!7XJxaH"

$ws.Range("D4").Value = "biswajit.deb@capgemini.com"
$ws.Range("J4").Value = "biyD9U%uMV1g"

$ws.Range("D5").Value = "debanjan.das@capgemini.com"
$ws.Range("J5").Value = "NTFJ}kI)1VgE"

$ws.Range("D6").Value = "dhiraj.kajari@capgemini.com"
$ws.Range("J6").Value = "d+h$%RNs6KHS"

$ws.Range("D7").Value = "mayur.bhorkar@capgemini.com"
$ws.Range("J7").Value = "m=A#UPx5}lgY"

$ws.Range("D8").Value = "manoj-kumar.b.s@capgemini.com"
$ws.Range("J8").Value = "{AnbwofBQ*d0"

# --- K / M helper columns for rows 2-8 ----------------------------------
$ws.Range("K2:K8").Value = 80
$ws.Range("M2:M8").Value = $true

# --- shared formulas for rows 2-7 (same fill order as the original author) --
$ws.Range("A2:A7").Formula = "=PROPER(IFERROR(LEFT(C2,FIND(CHAR(46),C2)-1),C2))"
$ws.Range("B2:B7").Formula = '=IFERROR(PROPER(RIGHT(C2,LEN(C2)-FIND("@",SUBSTITUTE(C2,".","@",((LEN(C2)-LEN(SUBSTITUTE(C2,".","")))/LEN("\")))))), "Unknown")'
$ws.Range("C2:C7").Formula = "=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D2,FIND(CHAR(64),D2)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))"
$ws.Range("E2:E7").Formula = "=LEFT(H2,FIND(CHAR(46),H2)-1)"
$ws.Range("F2:F7").Formula = '=CONCATENATE("ITPartner\",I2)'
$ws.Range("H2:H7").Formula = "=RIGHT(D2,LEN(D2)-FIND(CHAR(64),D2))"
$ws.Range("I2:I7").Formula = "=PROPER(E2)"
$ws.Range("P2:P7").Formula = "=COUNTIF(D:D,D2)"

# --- row 8 uses standalone (non-shared) formulas, entered individually -----
$ws.Range("A8").Formula = "=PROPER(IFERROR(LEFT(C8,FIND(CHAR(46),C8)-1),C8))"
$ws.Range("B8").Formula = '=IFERROR(PROPER(RIGHT(C8,LEN(C8)-FIND("@",SUBSTITUTE(C8,".","@",((LEN(C8)-LEN(SUBSTITUTE(C8,".","")))/LEN("\")))))), "Unknown")'
$ws.Range("C8").Formula = "=SUBSTITUTE(SUBSTITUTE(LOWER(LEFT(D8,FIND(CHAR(64),D8)-1)),CHAR(45),CHAR(46)),CHAR(95),CHAR(46))"
$ws.Range("E8").Formula = "=LEFT(H8,FIND(CHAR(46),H8)-1)"
$ws.Range("F8").Formula = '=CONCATENATE("ITPartner\",I8)'
$ws.Range("H8").Formula = "=RIGHT(D8,LEN(D8)-FIND(CHAR(64),D8))"
$ws.Range("I8").Formula = "=PROPER(E8)"
$ws.Range("P8").Formula = "=COUNTIF(D:D,D8)"

$excel.Calculate()

Write-Host "Email table populated (rows 2-8)"
